# Apply targeted cell edits to the "Risk Evaluation" sheet described by the diff.
# Summary of the change (5th data column, "N", plus a few scattered cells):
#   - Column N (rows 4-16) had its literal entries removed (now blank).
#   - Row 8's target-cost figures for columns H/J/L were recomputed to new values.
#   - Row 9's "Project life time" for column F dropped from 9 to 0.
#   - Row 11's "C- samples" placeholder text in F/J/L was removed.
#   - Row 12's "SOP Hella" text in F/J/L was removed, and H12's text was replaced.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk Evaluation")

# Clear out column N values (now blank) for rows 4-16.
$ws.Range("N4:N16").Value = $null

# Row 8: updated numeric target costs.
$ws.Range("H8").Value = 4.705436171336311
$ws.Range("J8").Value = 6.961506961506961
$ws.Range("L8").Value = 4.095004095004094

# Row 9: project life time changed.
$ws.Range("F9").Value = 0

# Row 11: clear text values that were removed.
$ws.Range("F11").Value = $null
$ws.Range("J11").Value = $null
$ws.Range("L11").Value = $null

# Row 12: clear text values that were removed, and update H12 to its new text value.
$ws.Range("F12").Value = $null
$ws.Range("J12").Value = $null
$ws.Range("L12").Value = $null
$ws.Range("H12").Value = "9-08-02 00:00:00"
